$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at C:D. This shifts the pre-existing column C
# (values, styles, and width formatting) two positions right to column E.
$ws.Range("C:D").Insert()

# Match column C's original custom width (~8.0 chars in the saved XML) on the
# two newly-inserted columns and re-assert it on the shifted-to column too.
$ws.Columns("C").ColumnWidth = 7.1667
$ws.Columns("D").ColumnWidth = 7.1667
$ws.Columns("E").ColumnWidth = 7.1667

# New header row: two new date headers ("Jun_17", "Jun_15") go in B1/C1, and
# D1 gets back the date header that used to sit in B1 ("Jun_13") before the
# new columns pushed things over (Insert() leaves the new D1 cell blank).
$ws.Range("D1").Value() = "Jun_13"
$ws.Range("C1").Value() = "Jun_15"
$ws.Range("B1").Value() = "Jun_17"

# Fill the two new columns (rows 2-27) with the same "UN" filler value used
# throughout column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value() = "UN"
    $ws.Cells.Item($r, 4).Value() = "UN"
}
